$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.911025
$ws.Range("H2").Value = 32.733075
$ws.Range("I2").Value = 0.03114956057965708
$ws.Range("J2").Value = 0.03274614149636444
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 1330.367519447475
$ws.Range("R2").Value = 11973.30767502727
$ws.Range("S2").Value = 0.007109052866026918
$ws.Range("T2").Value = 0.007923650762563055
$ws.Range("G3").Value = 10.911025
$ws.Range("H3").Value = 32.733075
$ws.Range("I3").Value = 0.03114956057965708
$ws.Range("J3").Value = 0.03274614149636444
$ws.Range("M3").Value = 147.91433
$ws.Range("N3").Value = 443.74299
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 1613.89695248825
$ws.Range("R3").Value = 14525.07257239425
$ws.Range("S3").Value = 0.008624142267336588
$ws.Range("T3").Value = 0.009612348190515638
$ws.Range("G4").Value = 10.911025
$ws.Range("H4").Value = 32.733075
$ws.Range("I4").Value = 0.03114956057965708
$ws.Range("J4").Value = 0.03274614149636444
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 911.1248080202334
$ws.Range("R4").Value = 8200.1232721821
$ws.Range("S4").Value = 0.004868755688243633
$ws.Range("T4").Value = 0.005426646903449656
$ws.Range("G5").Value = 10.911025
$ws.Range("H5").Value = 32.733075
$ws.Range("I5").Value = 0.03114956057965708
$ws.Range("J5").Value = 0.03274614149636444
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 993.6502710490627
$ws.Range("R5").Value = 5961.901626294375
$ws.Range("S5").Value = 0.00530974501704877
$ws.Range("T5").Value = 0.00394544495553571
$ws.Range("G6").Value = 10.911025
$ws.Range("H6").Value = 32.733075
$ws.Range("I6").Value = 0.03114956057965708
$ws.Range("J6").Value = 0.03274614149636444
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 980.19880482075
$ws.Range("R6").Value = 8821.789243386749
$ws.Range("S6").Value = 0.005237864741001166
$ws.Range("T6").Value = 0.005838050684300382
$ws.Range("G7").Value = 126.153142
$ws.Range("H7").Value = 378.459426
$ws.Range("I7").Value = 0.3601508510009905
$ws.Range("J7").Value = 0.3786105006764219
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 15381.69352494794
$ws.Range("R7").Value = 138435.2417245314
$ws.Range("S7").Value = 0.082194785154777
$ws.Range("T7").Value = 0.09161315639988227
$ws.Range("G8").Value = 126.153142
$ws.Range("H8").Value = 378.459426
$ws.Range("I8").Value = 0.3601508510009905
$ws.Range("J8").Value = 0.3786105006764219
$ws.Range("M8").Value = 147.91433
$ws.Range("N8").Value = 443.74299
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 18659.85747632486
$ws.Range("R8").Value = 167938.7172869237
$ws.Range("S8").Value = 0.09971223089302012
$ws.Range("T8").Value = 0.1111378560888241
$ws.Range("G9").Value = 126.153142
$ws.Range("H9").Value = 378.459426
$ws.Range("I9").Value = 0.3601508510009905
$ws.Range("J9").Value = 0.3786105006764219
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 10534.41425401365
$ws.Range("R9").Value = 94809.7282861228
$ws.Range("S9").Value = 0.05629249568233111
$ws.Range("T9").Value = 0.06274282731409238
$ws.Range("G10").Value = 126.153142
$ws.Range("H10").Value = 378.459426
$ws.Range("I10").Value = 0.3601508510009905
$ws.Range("J10").Value = 0.3786105006764219
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 11488.57268148418
$ws.Range("R10").Value = 68931.43608890506
$ws.Range("S10").Value = 0.06139120908618081
$ws.Range("T10").Value = 0.04561718791120725
$ws.Range("G11").Value = 126.153142
$ws.Range("H11").Value = 378.459426
$ws.Range("I11").Value = 0.3601508510009905
$ws.Range("J11").Value = 0.3786105006764219
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 11333.04698804946
$ws.Range("R11").Value = 101997.4228924451
$ws.Range("S11").Value = 0.0605601301846814
$ws.Range("T11").Value = 0.06749947296241583
$ws.Range("G12").Value = 48.19780633333333
$ws.Range("H12").Value = 144.593419
$ws.Range("I12").Value = 0.1375984830193998
$ws.Range("J12").Value = 0.1446511382757993
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 5876.697748789546
$ws.Range("R12").Value = 52890.27973910591
$ws.Range("S12").Value = 0.03140316819457326
$ws.Range("T12").Value = 0.03500153146995659
$ws.Range("G13").Value = 48.19780633333333
$ws.Range("H13").Value = 144.593419
$ws.Range("I13").Value = 0.1375984830193998
$ws.Range("J13").Value = 0.1446511382757993
$ws.Range("M13").Value = 147.91433
$ws.Range("N13").Value = 443.74299
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 7129.146231264755
$ws.Range("R13").Value = 64162.3160813828
$ws.Range("S13").Value = 0.03809584697974783
$ws.Range("T13").Value = 0.04246109751329868
$ws.Range("G14").Value = 48.19780633333333
$ws.Range("H14").Value = 144.593419
$ws.Range("I14").Value = 0.1375984830193998
$ws.Range("J14").Value = 0.1446511382757993
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 4024.756339798938
$ws.Range("R14").Value = 36222.80705819045
$ws.Range("S14").Value = 0.02150699349935333
$ws.Range("T14").Value = 0.02397139375006927
$ws.Range("G15").Value = 48.19780633333333
$ws.Range("H15").Value = 144.593419
$ws.Range("I15").Value = 0.1375984830193998
$ws.Range("J15").Value = 0.1446511382757993
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 4389.300118649429
$ws.Range("R15").Value = 26335.80071189658
$ws.Range("S15").Value = 0.02345499730878614
$ws.Range("T15").Value = 0.0174284076762483
$ws.Range("G16").Value = 48.19780633333333
$ws.Range("H16").Value = 144.593419
$ws.Range("I16").Value = 0.1375984830193998
$ws.Range("J16").Value = 0.1446511382757993
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 4329.880296572989
$ws.Range("R16").Value = 38968.92266915691
$ws.Range("S16").Value = 0.02313747703693918
$ws.Range("T16").Value = 0.02578870786622649
$ws.Range("G17").Value = 51.234875
$ws.Range("H17").Value = 102.46975
$ws.Range("I17").Value = 0.1462689199780642
$ws.Range("J17").Value = 0.102510654211286
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 6247.003701572625
$ws.Range("R17").Value = 37482.02220943575
$ws.Range("S17").Value = 0.0333819631940428
$ws.Range("T17").Value = 0.02480471244229715
$ws.Range("G18").Value = 51.234875
$ws.Range("H18").Value = 102.46975
$ws.Range("I18").Value = 0.1462689199780642
$ws.Range("J18").Value = 0.102510654211286
$ws.Range("M18").Value = 147.91433
$ws.Range("N18").Value = 443.74299
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 7578.372208258749
$ws.Range("R18").Value = 45470.2332495525
$ws.Range("S18").Value = 0.04049636501146379
$ws.Range("T18").Value = 0.03009112086154722
$ws.Range("G19").Value = 51.234875
$ws.Range("H19").Value = 102.46975
$ws.Range("I19").Value = 0.1462689199780642
$ws.Range("J19").Value = 0.102510654211286
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 4278.366665672167
$ws.Range("R19").Value = 25670.199994033
$ws.Range("S19").Value = 0.02286220488842263
$ws.Range("T19").Value = 0.01698792892310791
$ws.Range("G20").Value = 51.234875
$ws.Range("H20").Value = 102.46975
$ws.Range("I20").Value = 0.1462689199780642
$ws.Range("J20").Value = 0.102510654211286
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 4665.881292629688
$ws.Range("R20").Value = 18663.52517051875
$ws.Range("S20").Value = 0.02493295746553294
$ws.Range("T20").Value = 0.01235107786947928
$ws.Range("G21").Value = 51.234875
$ws.Range("H21").Value = 102.46975
$ws.Range("I21").Value = 0.1462689199780642
$ws.Range("J21").Value = 0.102510654211286
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 4602.71727359625
$ws.Range("R21").Value = 27616.3036415775
$ws.Range("S21").Value = 0.02459542941860202
$ws.Range("T21").Value = 0.01827581411485444
$ws.Range("G22").Value = 113.7817686666667
$ws.Range("H22").Value = 341.3453060000001
$ws.Range("I22").Value = 0.3248321854218885
$ws.Range("J22").Value = 0.3414815653401283
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 13873.26757471638
$ws.Range("R22").Value = 124859.4081724474
$ws.Range("S22").Value = 0.0741342457414751
$ws.Range("T22").Value = 0.08262899205724546
$ws.Range("G23").Value = 113.7817686666667
$ws.Range("H23").Value = 341.3453060000001
$ws.Range("I23").Value = 0.3248321854218885
$ws.Range("J23").Value = 0.3414815653401283
$ws.Range("M23").Value = 147.91433
$ws.Range("N23").Value = 443.74299
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 16829.95407854499
$ws.Range("R23").Value = 151469.586706905
$ws.Range("S23").Value = 0.08993382018742638
$ws.Range("T23").Value = 0.1002389764624957
$ws.Range("G24").Value = 113.7817686666667
$ws.Range("H24").Value = 341.3453060000001
$ws.Range("I24").Value = 0.3248321854218885
$ws.Range("J24").Value = 0.3414815653401283
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 9501.343103202429
$ws.Range("R24").Value = 85512.08792882186
$ws.Range("S24").Value = 0.05077209825971937
$ws.Range("T24").Value = 0.05658986965972417
$ws.Range("G25").Value = 113.7817686666667
$ws.Range("H25").Value = 341.3453060000001
$ws.Range("I25").Value = 0.3248321854218885
$ws.Range("J25").Value = 0.3414815653401283
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 10361.93073300401
$ws.Range("R25").Value = 62171.58439802407
$ws.Range("S25").Value = 0.05537079964612208
$ws.Range("T25").Value = 0.0411436785469588
$ws.Range("G26").Value = 113.7817686666667
$ws.Range("H26").Value = 341.3453060000001
$ws.Range("I26").Value = 0.3248321854218885
$ws.Range("J26").Value = 0.3414815653401283
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 10221.65687068426
$ws.Range("R26").Value = 91994.91183615835
$ws.Range("S26").Value = 0.05462122158714553
$ws.Range("T26").Value = 0.06088004861370414
